# Fruta / hortaliza, semanal
# Insert a new weekly record as row 162 in the Chirimoya (Vega Modelo de Temuco)
# price sheet, pushing the existing rows 162-191 down to 163-192.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(162).Insert()

$ws.Range('A162').Value = 10
$ws.Range('B162').Value = 'Vega Modelo de Temuco'
$ws.Range('C162').Value = 'La Araucanía'
$ws.Range('D162').Value = 45154
$ws.Range('E162').Value = 9
$ws.Range('F162').Value = 'Fruta'
$ws.Range('G162').Value = 100107
$ws.Range('H162').Value = 'Otros'
$ws.Range('I162').Value = 100107002
$ws.Range('J162').Value = 'Chirimoya'
$ws.Range('K162').Value = 'Cultivar IV Región'
$ws.Range('L162').Value = 'Especial'
$ws.Range('M162').Value = 50
$ws.Range('N162').Value = 3500
$ws.Range('O162').Value = 3500
$ws.Range('P162').Value = 3500
$ws.Range('Q162').Value = '$/kilo (en caja de 15 kilos)'
$ws.Range('R162').Value = 'Provincia del Elquí'
$ws.Range('S162').Value = 3500
$ws.Range('T162').Value = 1
